# POE2 to provider CPR
# Rename the Pharmacy "Alien" positive test case into a Provider-owned
# "POE2" case: clear the old Pharmacy cell, rename the Provider cell with
# a _POE2 suffix, and add a new Pharmacy row with the renamed "_POE2" case.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 9: clear the Pharmacy test case name (E9), rename the Provider
# test case name (G9) to carry the _POE2 suffix.
$ws2.Cells.Item(9, 5).Value = $null
$ws2.Cells.Item(9, 7).Value = "Prov_Patient_App_Create_PCase_Positive_Alien_POE2"

# New row 10: the Pharmacy module/test-case pair, renamed with the _POE2
# suffix, now on its own row.
$ws2.Cells.Item(10, 4).Value = "CPR\Pharmacy"
$ws2.Cells.Item(10, 5).Value = "Phar_Patient_App_Create_PCase_Positive_Alien_POE2"

# Widen column E to fit the longer renamed test case text.
$ws2.Columns.Item(5).ColumnWidth = 55.1428571428571

# Sheet1: move the active selection from row 3 to A2:B2 (last edited cell).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A2:B2").Select()
